$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header labels in row 1 to short variable-style names
$ws.Range("A1").Value = "R_p"
$ws.Range("B1").Value = "R_sp"
$ws.Range("C1").Value = "Q_p"
$ws.Range("D1").Value = "Q_sp"
$ws.Range("E1").Value = "A_p"
$ws.Range("F1").Value = "A_sp"

# Update A2 value (new total for R_p)
$ws.Range("A2").Value = 1828946000

# Move the active selection to A2 (matching the saved selection state)
$ws.Range("A2").Select()
